# Natmi following Dr Hou advice
#
# The "Cd177 -> Pecam1" LR-pair sheet gains a new sending cluster ("M1"),
# which duplicates the full set of target clusters that "M2" already had
# (ECs, FAPs, M1, M2, Neutro, sCs), and every numeric column is recomputed
# for both sending clusters ("M1" rows first, then the refreshed "M2" rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..D are the categorical labels (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster); E..T are the recomputed metrics.
$data = @(
  @("M1","Cd177","Pecam1","ECs",   2,0.6666666666666666,3.926606666666667,11.77982,         0.6634356070430004,0.6634356070430004,2,1,442.844055,       885.68811,         0.5956080271087346,0.496306158069548, 1738.8744186567,  10433.2465119402, 0.3951475730245672, 0.3292671772580499),
  @("M1","Cd177","Pecam1","FAPs",  2,0.6666666666666666,3.926606666666667,11.77982,         0.6634356070430004,0.6634356070430004,3,1,2.568168333333333, 7.704505,          0.003454086505236688,0.004317313547742695,10.08418689878889,90.75768208910002,0.002291563977380738,0.002864259534341644),
  @("M1","Cd177","Pecam1","M1",    2,0.6666666666666666,3.926606666666667,11.77982,         0.6634356070430004,0.6634356070430004,3,1,3.891101666666666, 11.673305,         0.005233380375768717,0.006541279137781407,15.27882574501111,137.5094317051,  0.003472010886485044,0.004339717495611722),
  @("M1","Cd177","Pecam1","M2",    2,0.6666666666666666,3.926606666666667,11.77982,         0.6634356070430004,0.6634356070430004,3,1,6.062445333333333, 18.187336,         0.008153753141026635,0.01019149602864148, 23.80483826216889,214.24354435952, 0.005409490164795776,0.006761401354438089),
  @("M1","Cd177","Pecam1","Neutro",2,0.6666666666666666,3.926606666666667,11.77982,         0.6634356070430004,0.6634356070430004,3,1,285.0064136666667, 855.019241,        0.3833225394825256,0.4791204824644771, 1119.108083946291,10071.97275551662,0.2543098216748539,0.3178655881305556),
  @("M1","Cd177","Pecam1","sCs",   2,0.6666666666666666,3.926606666666667,11.77982,         0.6634356070430004,0.6634356070430004,2,1,3.143744,          6.287488,          0.004228213386707702,0.003523270751809444,12.34424614869333,74.06547689216,  0.002805147314917765,0.002337463270003547),
  @("M2","Cd177","Pecam1","ECs",   3,1,                  1.991988333333333,5.975965,         0.3365643929569996,0.3365643929569996,2,1,442.844055,       885.68811,         0.5956080271087346,0.496306158069548, 882.140191046025, 5292.84114627615, 0.2004604540841675, 0.1670389808114981),
  @("M2","Cd177","Pecam1","FAPs",  3,1,                  1.991988333333333,5.975965,         0.3365643929569996,0.3365643929569996,3,1,2.568168333333333, 7.704505,          0.003454086505236688,0.004317313547742695,5.115761358036111,46.041852222325, 0.00116252252785595, 0.001453054013401051),
  @("M2","Cd177","Pecam1","M1",    3,1,                  1.991988333333333,5.975965,         0.3365643929569996,0.3365643929569996,3,1,3.891101666666666, 11.673305,         0.005233380375768717,0.006541279137781407,7.751029123813888,69.75926211432498,0.001761369489283673,0.002201561642169685),
  @("M2","Cd177","Pecam1","M2",    3,1,                  1.991988333333333,5.975965,         0.3365643929569996,0.3365643929569996,3,1,6.062445333333333, 18.187336,         0.008153753141026635,0.01019149602864148, 12.07632037547111,108.68688337924, 0.002744262976230858,0.003430094674203393),
  @("M2","Cd177","Pecam1","Neutro",3,1,                  1.991988333333333,5.975965,         0.3365643929569996,0.3365643929569996,3,1,285.0064136666667, 855.019241,        0.3833225394825256,0.4791204824644771, 567.7294509491738,5109.565058542565,0.1290127178076718,0.1612548943339215),
  @("M2","Cd177","Pecam1","sCs",   3,1,                  1.991988333333333,5.975965,         0.3365643929569996,0.3365643929569996,2,1,3.143744,          6.287488,          0.004228213386707702,0.003523270751809444,6.262301370986666,37.57380822592,  0.001423066071789937,0.001185807481805897)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $cols.Count; $j++) {
    $ws.Range($cols[$j] + $r).Value = $row[$j]
  }
}
